$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "Guía de Aprendizaje Ciencia y Medición.docx"
$ws.Range("B2").Value = 24777
$ws.Range("C2").Value = "C:\Users\gjoan\Desktop\QUIMICA"
$ws.Range("D2").Value = "20/01/2021 05:52:53"
$ws.Range("E2").Value = "21/01/2021 12:04:32"
$ws.Range("F2").Value = ".docx"

# Row 3
$ws.Range("A3").Value = "Material de apoyo unidad I.pdf"
$ws.Range("B3").Value = 259760
$ws.Range("C3").Value = "C:\Users\gjoan\Desktop\QUIMICA"
$ws.Range("D3").Value = "21/01/2021 12:04:07"
$ws.Range("E3").Value = "21/01/2021 12:04:32"
$ws.Range("F3").Value = ".pdf"

# Row 4
$ws.Range("A4").Value = "Medición y Cifras Significativas.pdf"
$ws.Range("B4").Value = 339523
$ws.Range("C4").Value = "C:\Users\gjoan\Desktop\QUIMICA"
$ws.Range("D4").Value = "21/01/2021 12:05:47"
$ws.Range("E4").Value = "21/01/2021 12:12:28"
$ws.Range("F4").Value = ".pdf"

# Row 5
$ws.Range("A5").Value = "NORMA COGUANOR_NGO_4_010_2a._Revision.pdf"
$ws.Range("B5").Value = 503289
$ws.Range("C5").Value = "C:\Users\gjoan\Desktop\QUIMICA"
$ws.Range("D5").Value = "21/01/2021 12:06:18"
$ws.Range("E5").Value = "21/01/2021 09:07:02"
$ws.Range("F5").Value = ".pdf"

# Row 6
$ws.Range("A6").Value = "Nueva imagen de mapa de bits.bmp"
$ws.Range("B6").Value = 0
$ws.Range("C6").Value = "C:\Users\gjoan\Desktop\QUIMICA"
$ws.Range("D6").Value = "28/01/2021 10:52:33"
$ws.Range("E6").Value = "28/01/2021 10:53:33"
$ws.Range("F6").Value = ".bmp"

# Row 7
$ws.Range("A7").Value = "Nuevo documento de texto.txt"
$ws.Range("B7").Value = 0
$ws.Range("C7").Value = "C:\Users\gjoan\Desktop\QUIMICA"
$ws.Range("D7").Value = "28/01/2021 10:52:25"
$ws.Range("E7").Value = "28/01/2021 10:52:36"
$ws.Range("F7").Value = ".txt"

# Row 8
$ws.Range("A8").Value = "TABLA DE CONVERSIÓN DE UNIDADES.pdf"
$ws.Range("B8").Value = 26028
$ws.Range("C8").Value = "C:\Users\gjoan\Desktop\QUIMICA"
$ws.Range("D8").Value = "21/01/2021 08:52:56"
$ws.Range("E8").Value = "23/01/2021 07:22:04"
$ws.Range("F8").Value = ".pdf"

# Row 9
$ws.Range("A9").Value = "Tarea Preparatoria Primer Parcial.pdf"
$ws.Range("B9").Value = 473194
$ws.Range("C9").Value = "C:\Users\gjoan\Desktop\QUIMICA"
$ws.Range("D9").Value = "23/01/2021 07:05:59"
$ws.Range("E9").Value = "28/01/2021 10:52:25"
$ws.Range("F9").Value = ".pdf"

# Row 10
$ws.Range("A10").Value = "ExcelReadWrite.suo"
$ws.Range("B10").Value = 20992
$ws.Range("C10").Value = "F:\Descargas\ExcelReadWrite\ExcelReadWrite"
$ws.Range("D10").Value = 43620.19122685185
$ws.Range("E10").Value = "24/11/2020 08:29:44"
$ws.Range("F10").Value = ".suo"

# Row 11
$ws.Range("A11").Value = "ExcelReadWrite.xlsx"
$ws.Range("B11").Value = 9477
$ws.Range("C11").Value = "F:\Descargas\ExcelReadWrite\ExcelReadWrite"
$ws.Range("D11").Value = 43620.19122685185
$ws.Range("E11").Value = "24/11/2020 08:29:44"
$ws.Range("F11").Value = ".xlsx"

# Row 12
$ws.Range("A12").Value = "ExcelWriteStep.cs"
$ws.Range("B12").Value = 6238
$ws.Range("C12").Value = "F:\Descargas\ExcelReadWrite\ExcelReadWrite"
$ws.Range("D12").Value = 43620.19122685185
$ws.Range("E12").Value = "24/11/2020 08:29:44"
$ws.Range("F12").Value = ".cs"

# Apply the date/time number format to column D for rows 2-12 (maps to
# built-in numFmtId 22: "m/d/yyyy h:mm")
$ws.Range("D2:D12").NumberFormat = "m/d/yy h:mm"
